$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.918.27"
$ws.Range("E2").Value = "  +1.73%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.386.88"
$ws.Range("E3").Value = "  +8.93%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5: Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.98"
$ws.Range("E5").Value = "  +8.98%  "

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "633.03"
$ws.Range("E6").Value = "  +3.69%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  +26.78%  "

# Row 8: Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +3.61%  "

# Row 9: USDC
$ws.Range("E9").Value = "  +0.05%  "

# Row 10: Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.879"
$ws.Range("E10").Value = "  +12.02%  "

# Row 11: LidoStakedEther
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.382.35"
$ws.Range("E11").Value = "  +8.85%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +1.50%  "

# Row 13: WrappedBTC
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "98.689.96"
$ws.Range("E13").Value = "  +1.94%  "

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.34"
$ws.Range("E14").Value = "  +6.90%  "

# Row 15: ShibaInu
$ws.Range("E15").Value = "  +3.68%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.005.38"
$ws.Range("E16").Value = "  +8.65%  "

# Row 17: Toncoin
$ws.Range("E17").Value = "  +3.70%  "

# Row 18: WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.385.26"
$ws.Range("E18").Value = "  +9.15%  "

# Row 19: SuiNetwork
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.60"
$ws.Range("E19").Value = "  +0.15%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.28"
$ws.Range("E20").Value = "  +5.30%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "499.22"
$ws.Range("E21").Value = "  -3.10%  "

# Row 22: Polkadot
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +9.83%  "

# Row 23: PEPE
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000213"
$ws.Range("E23").Value = "  +9.88%  "

# Row 24: Uniswap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.40"
$ws.Range("E24").Value = "  +6.34%  "

# Row 25: NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.81"
$ws.Range("E25").Value = "  +5.21%  "

# Row 26: Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.91"
$ws.Range("E26").Value = "  +5.40%  "

# Row 27: Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("E27").Value = "  +3.76%  "

# Row 28: WrappedeETH
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.559.63"
$ws.Range("E28").Value = "  +8.55%  "

# Row 29: Stellar
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.285"
$ws.Range("E29").Value = "  +20.62%  "

# Row 30: Cronos
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.203"
$ws.Range("E30").Value = "  +16.07%  "

# Row 31: Dai
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.10%  "

# Row 32: Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.134"
$ws.Range("E32").Value = "  +8.52%  "

# Row 33: Binance-PegBSC-USD
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.994"
$ws.Range("E33").Value = "  -0.28%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.64"
$ws.Range("E34").Value = "  +6.86%  "

# Row 35: EthereumClassic
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.10"
$ws.Range("E35").Value = "  +6.14%  "

# Row 36: RenderToken
$ws.Range("E36").Value = "  +0.90%  "

# Row 37: Kaspa
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +0.86%  "

# Row 38: PancakeSwap
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").Value = "  +5.58%  "

# Row 39: PolygonEcosystemToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.473"
$ws.Range("E39").Value = "  +8.71%  "

# Row 40: Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "503.55"
$ws.Range("E40").Value = "  +2.96%  "

# Row 41: WhiteBITCoin
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.86"
$ws.Range("E41").Value = "  +2.66%  "

# Row 42: Fetch.AI
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").Value = "  +2.96%  "

# Row 43: MantraDAO
$ws.Range("E43").Value = "  +5.02%  "

# Row 44: dogwifhat
$ws.Range("E44").Value = "  +5.10%  "

# Row 45: ARBITRUM
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.790"
$ws.Range("E45").Value = "  +14.20%  "

# Row 46: USDe
$ws.Range("E46").Value = "  +0.01%  "

# Row 47: Monero
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.61"
$ws.Range("E47").Value = "  -0.71%  "

# Row 48: Stacks
$ws.Range("E48").Value = "  +3.10%  "

# Row 49: Mantle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.842"
$ws.Range("E49").Value = "  +15.63%  "

# Row 50: Filecoin
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.71"
$ws.Range("E50").Value = "  +8.16%  "

# Row 51: OKB
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.47"
$ws.Range("E51").Value = "  +4.35%  "
